# Generate Report for Handoff
#
# Moves the localization status from "In Translation" to "Ready for
# handoff" and refreshes the associated handoff timestamps across the
# Overview / zh-cn / de-de sheets. Excel auto-grows the status columns
# to fit the new (longer) text, so their widths are refreshed as well.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet --------------------------------------------------
# E2 = zh-cn status, F2 = de-de status, G2 = Latest HO Xliff Generate Date
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-16 02:52:50"

# --- zh-cn sheet -------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-16 02:52:45"

# --- de-de sheet -------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-16 02:52:50"

# --- Resize the status columns to fit the new text (mirrors Excel's
#     own column auto-grow once the cell text got longer) --------------
$wsOverview.Columns.Item(5).ColumnWidth = 16.38265482584637
$wsOverview.Columns.Item(6).ColumnWidth = 16.38265482584637
$wsZhCn.Columns.Item(3).ColumnWidth = 16.38265482584637
$wsDeDe.Columns.Item(3).ColumnWidth = 16.38265482584637
